$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.0000006487708416141559,
    0.00000001182604016585146,
    46751564.37907436,
    843412.7080118682,
    116494.0962686534,
    24111650.68432828,
    11335933.31210752,
    13681675.36496993,
    960689.5850188414,
    60168433.01643151,
    22019306.84109979,
    401376.2493998191,
    31278054.07264849,
    3757594.980442458,
    32595011.41656321,
    107439536.1809415,
    53743950.33948664,
    73130376.27633405,
    21142893.90270996,
    385400.6630348664,
    27657845.60669194,
    15305410.29521715,
    18763263.71127747,
    1255195.668910339,
    15236.36738720788,
    277.7342647744208,
    176884907.5101677,
    12382544.41762704,
    225713.7008099194,
    4139.603774203179,
    50.2491567677146,
    0.9159606260303816,
    0.01669647656609077,
    0.0003043496868748282,
    0.000005547801150389279,
    0.0000001017470295257811
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}
